# Example7.xlsx - "sample num 7 is cleaned up"
# The End-Time column (E) on Aircraft_scheduling previously computed
# End = Start + Duration via a formula (e.g. =C4+0.93), which rolled past
# midnight (values > 1) for several rows. This clean-up replaces those
# formulas with their already-calculated time-of-day values (i.e. the
# fractional/MOD(...,1) part), removing the formula so E just holds the
# literal end-of-day time. Column F (Duration, =MOD(E-C,1)) is left as a
# formula and recalculates automatically from the new E values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aircraft_scheduling")

$ws.Range("E4").Value  = 0.2911111111111111
$ws.Range("E5").Value  = 0.24888888888888891
$ws.Range("E6").Value  = 0.21174768518518519
$ws.Range("E9").Value  = 0.27902777777777776
$ws.Range("E10").Value = 0.18375
$ws.Range("E11").Value = 0.44291666666666668
$ws.Range("E12").Value = 0.35569444444444448
$ws.Range("E13").Value = 0.58777777777777784
$ws.Range("E14").Value = 0.73208333333333331
$ws.Range("E15").Value = 0.59902777777777783
$ws.Range("E16").Value = 0.69638888888888895
$ws.Range("E17").Value = 0.64722222222222225
$ws.Range("E18").Value = 0.79541666666666666

# Move the active selection to mirror the saved cursor position.
$null = $ws.Range("E19").Select()
